$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                                      $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "WARNING: not found -> $old"
    }
}

# --- Title & byline ---------------------------------------------------
Replace-Text "Success: Formula or Fortune?" "The Symphony of Science: Unfolding the Secrets of Our World"
Replace-Text "John D Rockefeller IV" "Laraine Peterson"

# --- Author name / e-mail line -----------------------------------------
Replace-Text "rockefeller" "laraine"
Replace-Text "john.d4@gmail" "peterson860@yahoo"

# --- First body paragraph (three lines separated by manual breaks) ----
Replace-Text '"Glory belongs to the steadfast and humble foot-soldier who climbs day by day the dusty hill of progress," remarked the esteemed Winston Churchill, capturing the essence of his timeworn conviction. We, the denizens of the modern world, grapple with a profound question that has haunted generations since the dawn of civilization: Is success a deliberate pursuit, a testament to untiring effort, or a stroke of unpredictable luck?' `
             'From the intricate patterns of a snowflake to the boundless expanse of the universe, our world is a tapestry of wonders that beckons us to explore its mysteries. Science, the systematic study of natural phenomena, offers us a lens through which we can decipher the secrets hidden within the symphony of existence. Like a conductor orchestrating a harmonious ensemble, science guides us in unraveling the fundamental principles that govern the universe, unveiling the intricate mechanisms that drive life, and illuminating the complex interactions that shape our societies.'

Replace-Text 'As we embark on the arduous journey of exploring the elusive nature of success, we will navigate through captivating tales of individuals who have scaled the pinnacle of their endeavors, driven by an unwavering belief in the formula of perseverance, determination, and calculated risks. We will delve into the annals of history, discovering the lives of luminaries who seemingly defied the odds, propelled by a fortuitous sequence of events or a fortuitous twist of fate. In doing so, we will ponder: Is success truly a puzzle with missing pieces, or is it an enigma that defies logical explanation?' `
             "Through the prism of science, we unravel the enigmas that have captivated humanity for eons. We probe the depths of matter, seeking to comprehend the fundamental building blocks of reality, the very essence of existence. We investigate the intricacies of life, unraveling the secrets of DNA, the blueprint of heredity, and exploring the miraculous processes that sustain life on Earth. We delve into the annals of time, piecing together the puzzle of our planet's history, and deciphering the stories etched in the ancient rocks."

Replace-Text 'Moreover, our investigation will traverse various disciplines, from the concrete realm of science and technology to the ethereal domain of the arts and humanities. We will examine the arduous ascent of those who have meticulously crafted their path to prominence, mastering the intricacies of their chosen field through relentless application and dedication. We will discern the trajectories of individuals who have seemingly been anointed by fortune, finding themselves in the right place at the right time, as though destined for greatness. By continually holding the dual lenses of formula and fortune to our observations, we will diligently strive to understand their differential roles in shaping the lives and legacies of successful individuals.' `
             'But the pursuit of scientific knowledge is not merely an academic endeavor; it is a catalyst for progress and innovation. The discoveries made through scientific inquiry have transformed our lives, revolutionizing industries, improving healthcare, and shaping the very fabric of our societies. Science empowers us to tackle global challenges, from climate change to pandemics, by equipping us with the knowledge and tools necessary to address these complex issues effectively.'

# --- Summary paragraph ---------------------------------------------------
Replace-Text 'Our journey through the enigmatic landscape of success has unveiled a tapestry woven with threads of toil, talent, and happenstance' `
             'In this essay, we embark on a journey through the diverse realms of science, exploring the fundamental principles that govern our universe, unraveling the intricate mysteries of life, and delving into the profound impact that scientific discoveries have had on our world'

Replace-Text ' We have encountered individuals who have meticulously navigated the labyrinth of obstacles, driven by a relentless pursuit of ' `
             ' From the exploration of matter to the study of life and the examination of our past, science '

Replace-Text 'excellence' 'illuminates the hidden wonders of our existence, driving progress and innovation, and empowering us to address global challenges'

Replace-Text '. We have also marveled at the trajectories of those who have fortuitously found themselves perched upon the pinnacle of achievement, beneficiaries of a favorable cosmic alignment. Whether a product of formulaic endeavor or a stroke of serendipitous luck, success remains an elusive enigma, an intricate dance between individual agency and the whims of fortune. It is a journey that demands our curiosity, compels our contemplation, and ultimately, challenges us to recognize that in the grand symphony of life, both the protagonist and the conductor play their unique and vital parts' `
             ''

# --- Trailing empty paragraph -------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
